# Apply crypto price/volume updates per the Apr 29 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.314.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.278.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.277.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.818.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.278.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.338.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.724"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.105"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0723"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.079.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "426.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.256"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.14%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.113"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.42%  "
